$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the date values in column B (rows 2-5): 38165 -> 38174
$ws.Range("B2").Value = 38174
$ws.Range("B3").Value = 38174
$ws.Range("B4").Value = 38174
$ws.Range("B5").Value = 38174

# Update the lpuser-* strings in column E (rows 2-5)
$ws.Range("E2").Value = "lpuser-5-23-2012-55611"
$ws.Range("E3").Value = "lpuser-5-23-2012-55645"
$ws.Range("E4").Value = "lpuser-5-23-2012-55681"
$ws.Range("E5").Value = "lpuser-5-23-2012-55713"
